$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 180
$ws1.Range("F6").Value = 683
$ws1.Range("F8").Value = 489
$ws1.Range("F9").Value = 86
$ws1.Range("F11").Value = 419
$ws1.Range("F14").Value = 119
$ws1.Range("F15").Value = 203

# Sheet 3: 本地生活 (index 3)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 6228
$ws3.Range("F3").Value = 740
$ws3.Range("F5").Value = 1831

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 6228
$ws4.Range("F3").Value = 740
$ws4.Range("F5").Value = 1831
$ws4.Range("F12").Value = 180
$ws4.Range("F15").Value = 684
$ws4.Range("F19").Value = 489
$ws4.Range("F21").Value = 86
$ws4.Range("F24").Value = 419
$ws4.Range("F29").Value = 119
$ws4.Range("F35").Value = 203
